{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,style\");\nawait context.sync();\n\n// Locate the Title, Author and Abstract paragraphs by their style, then\n// collapse each one's many single-word/space runs into one run holding\n// the full text (matching the canonical OOXML after the edit).\nconst replacements = {\n  \"Title\": \"Questions: Trigonometric identities (degrees)\",\n  \"Author\": \"Dzhemma Ruseva\",\n  \"Abstract\": \"A selection of questions on trigonometric identities, where angles are measured in degrees.\"\n};\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const style = p.style;\n  if (Object.prototype.hasOwnProperty.call(replacements, style)) {\n    p.getRange().insertText(replacements[style], \"Replace\");\n    delete replacements[style];\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The Title, Author and Abstract paragraphs each held their text split\n# across one run per word/space. Collapse each into a single run with\n# the full text (style/formatting is left untouched).\n$replacements = @{\n    \"Title\"    = \"Questions: Trigonometric identities (degrees)\"\n    \"Author\"   = \"Dzhemma Ruseva\"\n    \"Abstract\" = \"A selection of questions on trigonometric identities, where angles are measured in degrees.\"\n}\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $styleName = $p.Range.Style.NameLocal\n    if ($replacements.ContainsKey($styleName)) {\n        $newText = $replacements[$styleName]\n        $r = $p.Range\n        $r.MoveEnd(1, -1) | Out-Null\n        $r.Find.Execute($r.Text, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n    }\n}\n"}
